$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update PriceAuditFolder (row 2) / ProcessedFolder (row 3) values ---
# Set B3 first, then B2, so the shared-string table reuses the two freed
# slots (formerly "Price Audit Folder" / "Processed Price Audit Data") in
# the same order they appear in the diff.
$ws.Range("B3").Value = "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Processed"
$ws.Range("B2").Value = "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Loaded"

# --- Add new row 4: TemplateFile ---
$ws.Range("A4").Value = "TemplateFile"
$ws.Range("B4").Value = "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Template\02 February 2019-DiCarlo Distributors Template.xlsx"

# --- Hyperlinks: add in the order B3, B2, B4 so relationship ids line up ---
$ws.Hyperlinks.Add($ws.Range("B3"), "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Processed") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Loaded") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Template\02 February 2019-DiCarlo Distributors Template.xlsx") | Out-Null

# --- Column B is now much wider than A/C to fit the long UNC paths ---
$ws.Columns.Item(2).ColumnWidth = 147.417

# --- Selection moves to A6 ---
$ws.Range("A6").Select() | Out-Null
